# feat: add 2022-Q3 data
#
# 1. "总计" sheet: insert the new 2022-Q3 summary row above the existing
#    2022-Q2 row (which shifts down and has its "age" counter bumped).
# 2. Duplicate the existing "2022-Q2" sheet (placing the copy right after
#    it) so the historical data is preserved unchanged under the name
#    "2022-Q2", then rename/refill the original sheet with the new
#    "2022-Q3" fund-holding data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: update the "总计" (summary) sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$oldDate  = $summary.Range("B2").Value()
$oldCount = $summary.Range("C2").Value()
$oldValue = $summary.Range("D2").Value()

# Push the existing row down to row 3, carrying A2's style/format with it,
# then bump its "periods old" counter from 0 to 1.
$summary.Range("A2").Copy($summary.Range("A3"))
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = $oldDate
$summary.Range("C3").Value = $oldCount
$summary.Range("D3").Value = $oldValue

# Overwrite row 2 with the new 2022-Q3 figures.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# ---------------------------------------------------------------------
# Step 2: duplicate the "2022-Q2" sheet so its data survives unchanged,
# then turn the original sheet into "2022-Q3" with fresh data.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)

$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

$q3 = $q2

# Match the "总计" sheet's page margins (this sheet was effectively
# rebuilt from that template), matching 0.75in/1in/0.5in.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Replace the header row text (keep the existing bold/bordered style).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Re-apply the "总计" sheet's header style (style index already used
# elsewhere in the workbook) onto the new header row + A2 flag cell so we
# don't introduce any extra style entries.
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: the new fund-holding figures for 2022-Q3.
$q3.Range("A2").Value = 0

$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "070031"
$q3.Range("B2").ClearFormats()

$q3.Range("C2").NumberFormat = "@"
$q3.Range("C2").Value = "嘉实全球房地产（QDII）"
$q3.Range("C2").ClearFormats()

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.38"
$q3.Range("D2").ClearFormats()

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "94.39"
$q3.Range("E2").ClearFormats()

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "2.49"
$q3.Range("F2").ClearFormats()

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0095"
$q3.Range("G2").ClearFormats()

$q3.Range("H2").Value = 8

Write-Host "Sheets now:" ($wb.Worksheets | ForEach-Object { $_.Name }) -join "|"
